# Commit: "Updated curve names in sheets"
#
# 1) Rename the four cost-curve option labels on the hidden
#    "Cost curve options" sheet (A1:A4), plus the column header
#    "Marginal costs" -> "Cost-coverage relationship" on
#    "Programs cost and coverage"!E1 (these share the sharedStrings
#    table with the dropdown list that the data validation on that
#    sheet points at).
$wb = $excel.ActiveWorkbook

$curveSheet = $wb.Worksheets.Item("Cost curve options")
$curveSheet.Range("A1").Value = "Linear (constant marginal cost) [default]"
$curveSheet.Range("A2").Value = "Curved with increasing marginal cost"
$curveSheet.Range("A3").Value = "Curved with decreasing marginal cost"
$curveSheet.Range("A4").Value = "S-shaped (decreasing then increasing marginal cost)"

$costSheet = $wb.Worksheets.Item("Programs cost and coverage")
$costSheet.Range("E1").Value = "Cost-coverage relationship"

# Every program row's "Cost curve" column was defaulted to the shared
# string "Constant (default)" (row E2:E38); since the rename above is
# really a rename of that shared-string slot (same dropdown list
# entry), every cell that held the old label now holds the new one.
$costSheet.Range("E2:E38").Value = "Linear (constant marginal cost) [default]"

# 2) Column E on "Programs cost and coverage" widens (and becomes a
#    best-fit width) to accommodate the longer header text.
$costSheet.Columns("E").ColumnWidth = 31.5

# 3) View-state: the cursor moves to E1 on "Programs cost and
#    coverage" (the edited header cell) and back to the default A1 on
#    "Cost curve options" (its stale A3 selection is cleared).
$curveSheet.Range("A1").Select()
$costSheet.Range("E1").Select()

# Restore the workbook's originally active tab/sheet so the rest of
# the view state (activeTab / tabSelected) is unaffected by the
# navigation above.
$wb.Worksheets.Item(1).Activate()
